# Applies the cryptos list refresh (prices / volume% / three row reorders)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.038.77'
$ws.Range("E2").Value = '  +4.97%  '

$ws.Range("D3").Value = '2.618.36'
$ws.Range("E3").Value = '  +5.54%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.43'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.38%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.12'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.46%  '

$ws.Range("E7").Value = '  +0.01%  '

$ws.Range("E8").Value = '  +2.07%  '

$ws.Range("D9").Value = '2.617.58'
$ws.Range("E9").Value = '  +5.52%  '

$ws.Range("E10").Value = '  +15.36%  '

$ws.Range("E11").Value = '  +0.57%  '

$ws.Range("E12").Value = '  +5.16%  '

$ws.Range("E13").Value = '  +2.18%  '

$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.81'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.41%  '

$ws.Range("B15").Value = 'ShibaInu'
$ws.Range("C15").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000184'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +8.47%  '

$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = '3.018.31'
$ws.Range("E16").Value = '  +2.90%  '

$ws.Range("D17").Value = '71.088.74'
$ws.Range("E17").Value = '  +5.14%  '

$ws.Range("D18").Value = '2.613.79'
$ws.Range("E18").Value = '  +5.46%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '384.34'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +10.59%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.91'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +6.93%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.50'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.77%  '

$ws.Range("E22").Value = '  +3.92%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '72.19'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +2.08%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.44'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +6.41%  '

$ws.Range("E25").Value = '  +0.19%  '

$ws.Range("E26").Value = '  +11.89%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.76'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +11.20%  '

$ws.Range("D28").Value = '2.750.62'
$ws.Range("E28").Value = '  +5.65%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.74%  '

$ws.Range("D30").Value = '0.0₃0952'
$ws.Range("E30").Value = '  +7.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '529.03'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +7.41%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.08'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +4.99%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.33'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +7.16%  '

$ws.Range("E34").Value = '  +4.95%  '

$ws.Range("E35").Value = '  +0.06%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '163.98'
$ws.Range("D36").Style = "Normal"

$ws.Range("E37").Value = '  +0.07%  '

$ws.Range("E38").Value = '  +5.39%  '

$ws.Range("E39").Value = '  +1.76%  '

$ws.Range("E40").Value = '  +7.29%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.84'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +6.96%  '

$ws.Range("B42").Value = 'USDe'
$ws.Range("C42").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.13%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.07'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +6.85%  '

$ws.Range("E44").Value = '  +9.90%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.333'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.84%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '40.13'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.86%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '154.15'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.73%  '

$ws.Range("E48").Value = '  +4.46%  '

$ws.Range("D49").Value = '0.0₆0273'
$ws.Range("E49").Value = '  +8.11%  '

$ws.Range("E50").Value = '  +5.21%  '

$ws.Range("E51").Value = '  +7.51%  '
